$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format (matches source data which stores these as text strings,
# even when the content looks numeric) for Price (D) column cells being updated
# that would otherwise be auto-coerced to a Number by Excel.
$textRefs = @("D5","D7","D9","D10","D11","D15","D18","D22","D25","D29","D33","D41","D42","D43","D44","D48")
foreach ($ref in $textRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply updated Price (D) and Volume(1h) (E) values per row
$ws.Range("D2").Value = "25.911.83"
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").Value = "1.584.39"
$ws.Range("E3").Value = "  -2.23%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "210.38"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("D7").Value = "0.479"
$ws.Range("E7").Value = "  -2.50%  "
$ws.Range("E8").Value = "  -0.66%  "
$ws.Range("D9").Value = "0.0614"
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").Value = "18.08"
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("D11").Value = "0.0791"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").Value = "1.804.36"
$ws.Range("E12").Value = "  -2.21%  "
$ws.Range("D13").Value = "1.581.92"
$ws.Range("E13").Value = "  -2.38%  "
$ws.Range("E14").Value = "  -2.77%  "
$ws.Range("D15").Value = "0.508"
$ws.Range("E15").Value = "  -2.84%  "
$ws.Range("D16").Value = "25.863.21"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("E17").Value = "  -1.68%  "
$ws.Range("D18").Value = "59.80"
$ws.Range("E18").Value = "  -3.17%  "
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("E21").Value = "  -1.60%  "
$ws.Range("D22").Value = "9.35"
$ws.Range("E22").Value = "  -1.80%  "
$ws.Range("E23").Value = "  -1.56%  "
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").Value = "141.04"
$ws.Range("E25").Value = "  -2.31%  "
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("E28").Value = "  -0.86%  "
$ws.Range("D29").Value = "6.45"
$ws.Range("E29").Value = "  -3.13%  "
$ws.Range("E30").Value = "  -5.59%  "
$ws.Range("E31").Value = "  -1.41%  "
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").Value = "3.02"
$ws.Range("E33").Value = "  -2.39%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").Value = "  -2.45%  "
$ws.Range("D36").Value = "1.096.68"
$ws.Range("E36").Value = "  -2.90%  "
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("E38").Value = "  -2.25%  "
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("E40").Value = "  -3.41%  "
$ws.Range("D41").Value = "0.774"
$ws.Range("E41").Value = "  -8.37%  "
$ws.Range("D42").Value = "0.811"
$ws.Range("E42").Value = "  +7.34%  "
$ws.Range("D43").Value = "93.76"
$ws.Range("E43").Value = "  -4.17%  "
$ws.Range("D44").Value = "5.13"
$ws.Range("D45").Value = "1.716.71"
$ws.Range("E45").Value = "  -2.25%  "
$ws.Range("D46").Value = "0.0₆0113"
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").Value = "53.09"
$ws.Range("E48").Value = "  -1.91%  "
$ws.Range("E49").Value = "  -1.25%  "
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("E51").Value = "  -0.35%  "
